$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A3").Value = "Kun jij dit even regelen?"
$ws.Range("B3").Value = "mailmind.test@zohomail.eu"
$ws.Range("C3").Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Range("D3").Value = "Overig"
$ws.Range("E3").Value = "Beste klant,`nDank je wel voor je e-mail. Kun je alsjeblieft meer details geven over wat je precies geregeld wilt hebben? Op die manier kan ik je beter helpen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$ws.Range("F3").Value = "2025-08-01 22:59:31"
$ws.Range("G3").Value = "Ja"
$ws.Range("H3").Value = "Nee"
$ws.Range("I3").Value = "Ja"
$ws.Range("J3").Value = "Ja"

# Extend conditional formatting ranges to include the new row
$colsToExtend = "D", "G", "H", "I", "J"
foreach ($col in $colsToExtend) {
    $srcCell = $ws.Range($col + "2")
    $newRange = $ws.Range($col + "2:" + $col + "3")
    $fcs = $srcCell.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 2
